$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 43406.168
$ws.Range("J87").Value = 48023.2
$ws.Range("L87").Value = 48023.2
$ws.Range("N87").Value = -50519.2
$ws.Range("H90").Value = 43406.168
$ws.Range("J90").Value = 48023.2
$ws.Range("L90").Value = 144069.6
$ws.Range("N90").Value = -156549.6
$ws.Range("H92").Value = 325.7
$ws.Range("I92").Value = 206.33333
$ws.Range("J92").Value = 1400
$ws.Range("K92").Value = 206.33333
$ws.Range("L92").Value = 1400
$ws.Range("M92").Value = 1041.66667
$ws.Range("N92").Value = -3896
$ws.Range("H100").Value = 3084.3333
$ws.Range("I100").Value = 1750
$ws.Range("K100").Value = 1750
$ws.Range("M100").Value = -1209
$ws.Range("H101").Value = 5242
$ws.Range("I101").Value = 6656
$ws.Range("J101").Value = 1000
$ws.Range("K101").Value = 19968
$ws.Range("L101").Value = 3000
$ws.Range("M101").Value = -18346
$ws.Range("N101").Value = -6244
$ws.Range("H103").Value = 639
$ws.Range("I103").Value = 626.9
$ws.Range("K103").Value = 1880.7
$ws.Range("M103").Value = -1294.7
$ws.Range("H136").Value = 35000
$ws.Range("J136").Value = 35000
$ws.Range("L136").Value = 35000
$ws.Range("N136").Value = -45200
$ws.Range("H137").Value = 1806.6285
$ws.Range("I137").Value = 2266.7693
$ws.Range("J137").Value = 1534.7273
$ws.Range("K137").Value = 6800.3079
$ws.Range("L137").Value = 4604.1819
$ws.Range("M137").Value = -4250.3079
$ws.Range("N137").Value = -9704.1819
$ws.Range("H138").Value = 2989159
$ws.Range("I138").Value = 9525547
$ws.Range("J138").Value = 5155.5435
$ws.Range("K138").Value = 28576641
$ws.Range("L138").Value = 15466.6305
$ws.Range("M138").Value = -28571501
$ws.Range("N138").Value = -25746.6305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 23333.334
$ws.Range("I3").Value = 23333.334
$ws.Range("K3").Value = 23333.334
$ws.Range("M3").Value = -23218.334
$ws.Range("H32").Value = 17758.041
$ws.Range("I32").Value = 18195.023
$ws.Range("K32").Value = 18195.023
$ws.Range("M32").Value = -17908.023

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2061.111
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("M99").Value = 498
$ws.Range("H134").Value = 2693.25
$ws.Range("I134").Value = 2713.9443
$ws.Range("J134").Value = 2507
$ws.Range("K134").Value = 8141.8329
$ws.Range("L134").Value = 7521
$ws.Range("M134").Value = -5606.8329
$ws.Range("N134").Value = -12591

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 1600
$ws.Range("I13").Value = 1600
$ws.Range("K13").Value = 1600
$ws.Range("M13").Value = -1461
$ws.Range("H31").Value = 24393060
$ws.Range("I31").Value = 33335440
$ws.Range("J31").Value = 4748.091
$ws.Range("K31").Value = 33335440
$ws.Range("L31").Value = 4748.091
$ws.Range("M31").Value = -33335145
$ws.Range("N31").Value = -5338.091
$ws.Range("H34").Value = 24393060
$ws.Range("I34").Value = 33335440
$ws.Range("J34").Value = 4748.091
$ws.Range("K34").Value = 33335440
$ws.Range("L34").Value = 4748.091
$ws.Range("M34").Value = -33335238
$ws.Range("N34").Value = -5152.091
$ws.Range("H41").Value = 1950
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = $null
$ws.Range("H58").Value = 1595.0714
$ws.Range("I58").Value = 1360.0769
$ws.Range("J58").Value = 4650
$ws.Range("K58").Value = 1360.0769
$ws.Range("L58").Value = 4650
$ws.Range("M58").Value = -1157.0769
$ws.Range("N58").Value = -5056
$ws.Range("H132").Value = 2331.7307
$ws.Range("I132").Value = 2171.7778
$ws.Range("J132").Value = 2691.625
$ws.Range("K132").Value = 6515.3334
$ws.Range("L132").Value = 8074.875
$ws.Range("M132").Value = -3985.3334
$ws.Range("N132").Value = -13134.875
$ws.Range("H134").Value = 1323.76
$ws.Range("I134").Value = 1273.6522
$ws.Range("K134").Value = 3820.9566
$ws.Range("M134").Value = -1285.9566
$ws.Range("H136").Value = 1595.0714
$ws.Range("I136").Value = 1360.0769
$ws.Range("J136").Value = 4650
$ws.Range("K136").Value = 4080.2307
$ws.Range("L136").Value = 13950
$ws.Range("M136").Value = -1530.2307
$ws.Range("N136").Value = -19050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1923.3334
$ws.Range("I13").Value = 1885
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 5655
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = -5487
$ws.Range("N13").Value = -6336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3284.1667
$ws.Range("I80").Value = 3843.5715
$ws.Range("J80").Value = 2501
$ws.Range("K80").Value = 3843.5715
$ws.Range("L80").Value = 2501
$ws.Range("M80").Value = -2845.5715
$ws.Range("N80").Value = -4497
$ws.Range("H83").Value = 3284.1667
$ws.Range("I83").Value = 3843.5715
$ws.Range("J83").Value = 2501
$ws.Range("K83").Value = 19217.8575
$ws.Range("L83").Value = 12505
$ws.Range("M83").Value = -14225.8575
$ws.Range("N83").Value = -22489
$ws.Range("H132").Value = 2061.4443
$ws.Range("I132").Value = 1454.5714
$ws.Range("J132").Value = 2911.0667
$ws.Range("K132").Value = 4363.7142
$ws.Range("L132").Value = 8733.2001
$ws.Range("M132").Value = -1833.7142
$ws.Range("N132").Value = -13793.2001
$ws.Range("H141").Value = 41214.285
$ws.Range("J141").Value = 41214.285
$ws.Range("L141").Value = 41214.285
$ws.Range("N141").Value = -51574.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 50000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = $null
$ws.Range("H136").Value = 4504.9707
$ws.Range("I136").Value = 4583.4062
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 13750.2186
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -11200.2186
$ws.Range("N136").Value = -14850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1886.174
$ws.Range("I132").Value = 2034.1
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 6102.299999999999
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -3572.299999999999
$ws.Range("N132").Value = -7760
$ws.Range("H140").Value = 43485.8
$ws.Range("J140").Value = 43485.8
$ws.Range("L140").Value = 43485.8
$ws.Range("N140").Value = -53845.8
$ws.Range("H141").Value = 70752.07000000001
$ws.Range("J141").Value = 70752.07000000001
$ws.Range("L141").Value = 70752.07000000001
$ws.Range("N141").Value = -81112.07000000001
